# Daily attendance processing - 2025-12-17 10:00:08
#
# This script applies the day's attendance-sync update to the
# "Session Analysis Results" sheet:
#   1. Six B1 sub-groups that were "Not Recorded" for session 11 (17/12/2025)
#      now have attendance recorded -> flip their row formatting from the
#      "Not Recorded" (pink) style to the "Recorded" (green) style, and fill
#      in Recorded By / Students / Status.
#   2. The "Recorded by" audit trail for a number of earlier rows now lists
#      the submitting user before "System" instead of after.
#   3. The roll-up statistics (overall Class Statistics block + the per-group
#      Group Statistics table) are refreshed to reflect the newly recorded
#      sessions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Newly-recorded sessions (row 168/195/222/249/276/303): copy the
#    formatting of an already-"Recorded" row onto them (this reuses the
#    existing green style instead of minting a new one), then fill in the
#    real values.
# ---------------------------------------------------------------------------
function Set-SessionRecorded {
    param($Row, $TemplateRow, $RecordedBy, $Students)

    $ws.Range("A${TemplateRow}:I${TemplateRow}").Copy()
    $ws.Range("A${Row}:I${Row}").PasteSpecial(-4122)

    $ws.Range("G$Row").Value = $RecordedBy
    $ws.Range("H$Row").Value = $Students
    $ws.Range("I$Row").Value = "Recorded"
}

# NOTE: positional args only - named args (e.g. "-Row 168") make this
# runtime's parameter binder pathologically slow.
Set-SessionRecorded 168 167 "dnasr281@gmail.com" "19/23"
Set-SessionRecorded 195 194 "dnasr281@gmail.com" "24/30"
Set-SessionRecorded 222 221 "dnasr281@gmail.com" "20/25"
Set-SessionRecorded 249 248 "dnasr281@gmail.com" "25/28"
Set-SessionRecorded 276 275 "dnasr281@gmail.com" "20/26"
Set-SessionRecorded 303 302 "dnasr281@gmail.com" "18/29"

# ---------------------------------------------------------------------------
# 2. "Recorded By" text reorder: "System, dnasr281@gmail.com" ->
#    "dnasr281@gmail.com, System" for every row where that applies.
# ---------------------------------------------------------------------------
$recordedByRows = 8,9,10,34,35,36,60,61,62,86,87,88,112,113,114,138,139,140,164,167,191,194,218,221,245,248,272,275,299,302
foreach ($r in $recordedByRows) {
    $ws.Range("G$r").Value = "dnasr281@gmail.com, System"
}

# ---------------------------------------------------------------------------
# 3a. Overall Class Statistics block (K/L columns).
# ---------------------------------------------------------------------------
$ws.Range("L6").Value = 129           # Recorded Sessions
$ws.Range("L7").Value = 3             # Missing Sessions
$ws.Range("L9").Value = "'40.6%"      # Coverage %
$ws.Range("L10").Value = "'71.5%"     # Average Attendance %

# ---------------------------------------------------------------------------
# 3b. Per-group Group Statistics table (rows 21-26: B1D1, B1D2, B1E1, B1E2,
#     B1F1, B1F2) - Recorded / Missing counts and Coverage / Avg Attendance
#     percentages.
# ---------------------------------------------------------------------------
$ws.Range("O21").Value = 11
$ws.Range("P21").Value = 0
$ws.Range("R21").Value = "'40.7%"
$ws.Range("S21").Value = "'77.1%"

$ws.Range("O22").Value = 11
$ws.Range("P22").Value = 0
$ws.Range("R22").Value = "'40.7%"
$ws.Range("S22").Value = "'74.2%"

$ws.Range("O23").Value = 11
$ws.Range("P23").Value = 0
$ws.Range("R23").Value = "'40.7%"

$ws.Range("O24").Value = 10
$ws.Range("P24").Value = 1
$ws.Range("R24").Value = "'37.0%"
$ws.Range("S24").Value = "'69.6%"

$ws.Range("O25").Value = 11
$ws.Range("P25").Value = 0
$ws.Range("R25").Value = "'40.7%"
$ws.Range("S25").Value = "'67.5%"

$ws.Range("O26").Value = 11
$ws.Range("P26").Value = 0
$ws.Range("R26").Value = "'40.7%"
$ws.Range("S26").Value = "'59.2%"
